$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 1399.9688
$ws.Range("I40").Value = 3033.3333
$ws.Range("J40").Value = 1023.03845
$ws.Range("K40").Value = 3033.3333
$ws.Range("L40").Value = 1023.03845
$ws.Range("M40").Value = -2858.3333
$ws.Range("N40").Value = -1373.03845

# Row 116
$ws.Range("H116").Value = 7188.8887
$ws.Range("I116").Value = 1800
$ws.Range("J116").Value = 9883.333000000001
$ws.Range("K116").Value = 1800
$ws.Range("L116").Value = 9883.333000000001
$ws.Range("M116").Value = 1642
$ws.Range("N116").Value = -16767.333

# Row 125
$ws.Range("H125").Value = 557.3
$ws.Range("I125").Value = 674.3333
$ws.Range("J125").Value = 381.75
$ws.Range("K125").Value = 6068.9997
$ws.Range("L125").Value = 3435.75
$ws.Range("M125").Value = -3608.9997
$ws.Range("N125").Value = -8355.75

# Row 137
$ws.Range("H137").Value = 1962729.2
$ws.Range("I137").Value = 2704311.2
$ws.Range("J137").Value = 2834
$ws.Range("K137").Value = 8112933.600000001
$ws.Range("L137").Value = 8502
$ws.Range("M137").Value = -8110383.600000001
$ws.Range("N137").Value = -13602

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 5542.222
$ws.Range("I45").Value = 5714.2856
$ws.Range("J45").Value = 4940
$ws.Range("K45").Value = 5714.2856
$ws.Range("L45").Value = 4940
$ws.Range("M45").Value = -5337.2856
$ws.Range("N45").Value = -5694

# Row 106
$ws.Range("H106").Value = 40303.125
$ws.Range("J106").Value = 40303.125
$ws.Range("L106").Value = 40303.125
$ws.Range("N106").Value = -42827.125

# Row 119
$ws.Range("H119").Value = 29349
$ws.Range("J119").Value = 29349
$ws.Range("L119").Value = 29349
$ws.Range("N119").Value = -39025

# Row 132
$ws.Range("H132").Value = 63462.266
$ws.Range("I132").Value = 43346.043
$ws.Range("J132").Value = 111741.2
$ws.Range("K132").Value = 130038.129
$ws.Range("L132").Value = 335223.6
$ws.Range("M132").Value = -127508.129
$ws.Range("N132").Value = -340283.6

# Row 135
$ws.Range("H135").Value = 50000
$ws.Range("J135").Value = 50000
$ws.Range("L135").Value = 50000
$ws.Range("N135").Value = -60140

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2839.2646
$ws.Range("I31").Value = 1322.6
$ws.Range("J31").Value = 7052.222
$ws.Range("K31").Value = 1322.6
$ws.Range("L31").Value = 7052.222
$ws.Range("M31").Value = -1027.6
$ws.Range("N31").Value = -7642.222

# Row 34
$ws.Range("H34").Value = 2839.2646
$ws.Range("I34").Value = 1322.6
$ws.Range("J34").Value = 7052.222
$ws.Range("K34").Value = 1322.6
$ws.Range("L34").Value = 7052.222
$ws.Range("M34").Value = -1120.6
$ws.Range("N34").Value = -7456.222

# Row 133
$ws.Range("H133").Value = 38115.6
$ws.Range("J133").Value = 38115.6
$ws.Range("L133").Value = 38115.6
$ws.Range("N133").Value = -43175.6

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 22316.805
$ws.Range("I5").Value = 30709.666
$ws.Range("K5").Value = 92128.99800000001
$ws.Range("M5").Value = -92016.99800000001

# Row 63
$ws.Range("H63").Value = 4733.3335
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 4733.3335
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 14200.0005
$ws.Range("M63").ClearContents()
$ws.Range("N63").Value = -15698.0005

# Row 66
$ws.Range("H66").Value = 4733.3335
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4733.3335
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 42600.0015
$ws.Range("M66").ClearContents()
$ws.Range("N66").Value = -50088.0015

# Row 114
$ws.Range("H114").Value = 14493255
$ws.Range("I114").Value = 457.86667
$ws.Range("J114").Value = 41667250
$ws.Range("K114").Value = 1373.60001
$ws.Range("L114").Value = 125001750
$ws.Range("M114").Value = 1880.39999
$ws.Range("N114").Value = -125008258

# Row 117
$ws.Range("H117").Value = 11112416
$ws.Range("I117").Value = 1276.3334
$ws.Range("J117").Value = 22223554
$ws.Range("K117").Value = 3829.0002
$ws.Range("L117").Value = 66670662
$ws.Range("M117").Value = -387.0001999999999
$ws.Range("N117").Value = -66677546

# Row 122
$ws.Range("H122").Value = 1191.375
$ws.Range("J122").Value = 1522.8823
$ws.Range("L122").Value = 13705.9407
$ws.Range("N122").Value = -18605.9407

# Row 126
$ws.Range("H126").Value = 3325
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 3325
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 9975
$ws.Range("M126").ClearContents()
$ws.Range("N126").Value = -19855

# Row 131
$ws.Range("H131").Value = 955.9231
$ws.Range("J131").Value = 962.16
$ws.Range("L131").Value = 2886.48
$ws.Range("N131").Value = -12966.48

# Row 132
$ws.Range("H132").Value = 1181.3572
$ws.Range("I132").Value = 613
$ws.Range("J132").Value = 1336.3636
$ws.Range("K132").Value = 5517
$ws.Range("L132").Value = 12027.2724
$ws.Range("M132").Value = -2987
$ws.Range("N132").Value = -17087.2724

# Row 135
$ws.Range("H135").Value = 22316.805
$ws.Range("I135").Value = 30709.666
$ws.Range("K135").Value = 276386.994
$ws.Range("M135").Value = -273851.994

$ws = $wb.Worksheets.Item("GSM")
# Row 100
$ws.Range("H100").Value = 33632.332
$ws.Range("J100").Value = 33632.332
$ws.Range("L100").Value = 33632.332
$ws.Range("N100").Value = -35796.332

# Row 101
$ws.Range("H101").Value = 42422
$ws.Range("J101").Value = 42422
$ws.Range("L101").Value = 42422
$ws.Range("N101").Value = -48912

# Row 102
$ws.Range("H102").Value = 1811.45
$ws.Range("I102").Value = 1336.2
$ws.Range("K102").Value = 1336.2
$ws.Range("M102").Value = 285.8

$ws = $wb.Worksheets.Item("LTW")
# Row 76
$ws.Range("H76").Value = 37370.824
$ws.Range("J76").Value = 37370.824
$ws.Range("L76").Value = 37370.824
$ws.Range("N76").Value = -38046.824

# Row 79
$ws.Range("H79").Value = 37370.824
$ws.Range("J79").Value = 37370.824
$ws.Range("L79").Value = 37370.824
$ws.Range("N79").Value = -39710.824

# Row 119
$ws.Range("H119").Value = 15000
$ws.Range("J119").Value = 15000
$ws.Range("L119").Value = 15000
$ws.Range("N119").Value = -24676

$ws = $wb.Worksheets.Item("WVR")
# Row 119
$ws.Range("H119").Value = 29060
$ws.Range("J119").Value = 29060
$ws.Range("L119").Value = 29060
$ws.Range("N119").Value = -38736

# Row 136
$ws.Range("H136").Value = 80599.48
$ws.Range("I136").Value = 42291.125
$ws.Range("K136").Value = 126873.375
$ws.Range("M136").Value = -124323.375
